# Generate Report for Handoff
#
# A new localized file "2291a6aa-fa1a-4622-b6e6-71faa470e61f.md" has been
# picked up by the handoff report. It needs to be inserted as a new row
# (row 7) in all three worksheets (Overview, zh-cn, de-de), pushing the
# existing "27a3f539-08b5-4bbe-9811-281f0ccb21ed.md" row and the
# ".localization-config" row down by one.

$wb = $excel.ActiveWorkbook

$newFileMd   = "2291a6aa-fa1a-4622-b6e6-71faa470e61f.md"
$newFileZh   = "2291a6aa-fa1a-4622-b6e6-71faa470e61f.ee29a6e0dc16ac4a8338fd020a8937d60ebc5a81.zh-cn.xlf"
$newFileDe   = "2291a6aa-fa1a-4622-b6e6-71faa470e61f.ee29a6e0dc16ac4a8338fd020a8937d60ebc5a81.de-de.xlf"
$zhDateTime  = "2016-03-04 10:17:48"
$deDateTime  = "2016-03-04 10:18:01"

$readyStatus = "Ready for handoff"
$includeStat = "Include"
$epochDate   = "0001-01-01 00:00:00"

function Add-Link($ws, $cellref, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellref), $url, [Type]::Missing, [Type]::Missing, $text) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new blank row at position 7; this shifts the existing rows 7-8
# (27a3f539...md, .localization-config) down to rows 8-9, carrying their
# styling along.
$ws1.Rows.Item(7).Insert()

$ws1.Range("A7").Value = $newFileMd
$ws1.Range("B7").Value = $readyStatus
$ws1.Range("C7").Value = $readyStatus

# Hyperlinks are not shifted automatically by Insert(), so rebuild the
# whole hyperlink collection for this sheet in the correct final order.
$ws1.Hyperlinks.Delete()

Add-Link $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/199a1f31-4b40-43bb-923a-18397fa6ca04.md" "199a1f31-4b40-43bb-923a-18397fa6ca04.md"
Add-Link $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/369213bbcf21725cdf2be78fee7674dca2193c34/e2e/1a074949-8159-4044-af6e-13f93f8e43e3.md" "1a074949-8159-4044-af6e-13f93f8e43e3.md"
Add-Link $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f808b105251aa7a0a36d71b489338eb774aeb18f/e2e/bea0c538-ddcb-4d77-977d-c7666a7b139c.md" "bea0c538-ddcb-4d77-977d-c7666a7b139c.md"
Add-Link $ws1 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/e74c1a4c-b419-43a7-8bd7-7f98d9726133.md" "e74c1a4c-b419-43a7-8bd7-7f98d9726133.md"
Add-Link $ws1 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/0a7d79782dc0e4eaac7da25566514f6648f5c4ee/e2e/0034157f-dd4d-417a-9df6-110b159bec4f.md" "0034157f-dd4d-417a-9df6-110b159bec4f.md"
Add-Link $ws1 "A7" "https://github.com/OpenLocalizationTest/oltest/blob/728ca3a2e4a4ae42a61684336ee819a2e04c0841/e2e/2291a6aa-fa1a-4622-b6e6-71faa470e61f.md" $newFileMd
Add-Link $ws1 "A8" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/e2e/27a3f539-08b5-4bbe-9811-281f0ccb21ed.md" "27a3f539-08b5-4bbe-9811-281f0ccb21ed.md"
Add-Link $ws1 "A9" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/.localization-config" ".localization-config"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(7).Insert()

$ws2.Range("A7").Value = $newFileMd
$ws2.Range("B7").Value = $readyStatus
$ws2.Range("C7").Value = $newFileZh
$ws2.Range("D7").Value = $zhDateTime
$ws2.Range("G7").Value = $epochDate
$ws2.Range("H7").Value = $includeStat

$ws2.Hyperlinks.Delete()

Add-Link $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/199a1f31-4b40-43bb-923a-18397fa6ca04.md" "199a1f31-4b40-43bb-923a-18397fa6ca04.md"
Add-Link $ws2 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/099d1bb2d354adf14802c5fb265e340f67b4bd7e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/199a1f31-4b40-43bb-923a-18397fa6ca04.9809858e68fcb9f6dd2bac79ef216efc8fdd31ee.zh-cn.xlf" "199a1f31-4b40-43bb-923a-18397fa6ca04.9809858e68fcb9f6dd2bac79ef216efc8fdd31ee.zh-cn.xlf"
Add-Link $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/369213bbcf21725cdf2be78fee7674dca2193c34/e2e/1a074949-8159-4044-af6e-13f93f8e43e3.md" "1a074949-8159-4044-af6e-13f93f8e43e3.md"
Add-Link $ws2 "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e8c1e35b1dbe525312190456d709acde1ad316c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.zh-cn.xlf" "1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.zh-cn.xlf"
Add-Link $ws2 "E3" "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3f7a838470b870e0b9033b02e7272599f7807bd6/e2e/1a074949-8159-4044-af6e-13f93f8e43e3.md" "1a074949-8159-4044-af6e-13f93f8e43e3.md"
Add-Link $ws2 "F3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/66b64a8b7952ca79a6706b3396245683e74c3395/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.zh-cn.xlf" "1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.zh-cn.xlf"
Add-Link $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f808b105251aa7a0a36d71b489338eb774aeb18f/e2e/bea0c538-ddcb-4d77-977d-c7666a7b139c.md" "bea0c538-ddcb-4d77-977d-c7666a7b139c.md"
Add-Link $ws2 "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e02a45df03dd737622128db4589d4834403e911d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/bea0c538-ddcb-4d77-977d-c7666a7b139c.61e8ee26ecc65f85dadc73a0d3c1f3f6c3ec4f66.zh-cn.xlf" "bea0c538-ddcb-4d77-977d-c7666a7b139c.61e8ee26ecc65f85dadc73a0d3c1f3f6c3ec4f66.zh-cn.xlf"
Add-Link $ws2 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/e74c1a4c-b419-43a7-8bd7-7f98d9726133.md" "e74c1a4c-b419-43a7-8bd7-7f98d9726133.md"
Add-Link $ws2 "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/099d1bb2d354adf14802c5fb265e340f67b4bd7e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/e74c1a4c-b419-43a7-8bd7-7f98d9726133.c2bc1a05f8e55b2535db8e200e39f576949ee6dd.zh-cn.xlf" "e74c1a4c-b419-43a7-8bd7-7f98d9726133.c2bc1a05f8e55b2535db8e200e39f576949ee6dd.zh-cn.xlf"
Add-Link $ws2 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/0a7d79782dc0e4eaac7da25566514f6648f5c4ee/e2e/0034157f-dd4d-417a-9df6-110b159bec4f.md" "0034157f-dd4d-417a-9df6-110b159bec4f.md"
Add-Link $ws2 "C6" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc3af3f65a0b196686e6acf4aad3b13f17bfbfa7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/0034157f-dd4d-417a-9df6-110b159bec4f.e05f1b9a6487db1a122b3f7a99a8da13ba3f092c.zh-cn.xlf" "0034157f-dd4d-417a-9df6-110b159bec4f.e05f1b9a6487db1a122b3f7a99a8da13ba3f092c.zh-cn.xlf"
Add-Link $ws2 "A7" "https://github.com/OpenLocalizationTest/oltest/blob/728ca3a2e4a4ae42a61684336ee819a2e04c0841/e2e/2291a6aa-fa1a-4622-b6e6-71faa470e61f.md" $newFileMd
Add-Link $ws2 "C7" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f523daf9b108236a3c1732329e900f2646e611f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/2291a6aa-fa1a-4622-b6e6-71faa470e61f.ee29a6e0dc16ac4a8338fd020a8937d60ebc5a81.zh-cn.xlf" $newFileZh
Add-Link $ws2 "A8" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/e2e/27a3f539-08b5-4bbe-9811-281f0ccb21ed.md" "27a3f539-08b5-4bbe-9811-281f0ccb21ed.md"
Add-Link $ws2 "C8" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b8f4befe585da0b687b4245c50d0728655fd3c88/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/ht/27a3f539-08b5-4bbe-9811-281f0ccb21ed.87a85293aeafa20d19b587592fc53056c03f5bcf.zh-cn.xlf" "27a3f539-08b5-4bbe-9811-281f0ccb21ed.87a85293aeafa20d19b587592fc53056c03f5bcf.zh-cn.xlf"
Add-Link $ws2 "A9" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/.localization-config" ".localization-config"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Rows.Item(7).Insert()

$ws3.Range("A7").Value = $newFileMd
$ws3.Range("B7").Value = $readyStatus
$ws3.Range("C7").Value = $newFileDe
$ws3.Range("D7").Value = $deDateTime
$ws3.Range("G7").Value = $epochDate
$ws3.Range("H7").Value = $includeStat

$ws3.Hyperlinks.Delete()

Add-Link $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/199a1f31-4b40-43bb-923a-18397fa6ca04.md" "199a1f31-4b40-43bb-923a-18397fa6ca04.md"
Add-Link $ws3 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97007c52653b3a1097edffc98ab2dcd10ec418bb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/199a1f31-4b40-43bb-923a-18397fa6ca04.9809858e68fcb9f6dd2bac79ef216efc8fdd31ee.de-de.xlf" "199a1f31-4b40-43bb-923a-18397fa6ca04.9809858e68fcb9f6dd2bac79ef216efc8fdd31ee.de-de.xlf"
Add-Link $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/369213bbcf21725cdf2be78fee7674dca2193c34/e2e/1a074949-8159-4044-af6e-13f93f8e43e3.md" "1a074949-8159-4044-af6e-13f93f8e43e3.md"
Add-Link $ws3 "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d90273a60bb4f19eefb2ad3290e93bc88f39f56f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.de-de.xlf" "1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.de-de.xlf"
Add-Link $ws3 "E3" "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4522e76346eb3d36cdcf7a47b52df48ee568231c/e2e/1a074949-8159-4044-af6e-13f93f8e43e3.md" "1a074949-8159-4044-af6e-13f93f8e43e3.md"
Add-Link $ws3 "F3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d9066770de28ab74054eb4d19d550e9bac12e7ea/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.de-de.xlf" "1a074949-8159-4044-af6e-13f93f8e43e3.31402265285413c8f21d97a5e66318f15fe231be.de-de.xlf"
Add-Link $ws3 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/f808b105251aa7a0a36d71b489338eb774aeb18f/e2e/bea0c538-ddcb-4d77-977d-c7666a7b139c.md" "bea0c538-ddcb-4d77-977d-c7666a7b139c.md"
Add-Link $ws3 "C4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b49c8b3c0ba97ea9bbb82fbfecf69482ab3e02c5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/bea0c538-ddcb-4d77-977d-c7666a7b139c.61e8ee26ecc65f85dadc73a0d3c1f3f6c3ec4f66.de-de.xlf" "bea0c538-ddcb-4d77-977d-c7666a7b139c.61e8ee26ecc65f85dadc73a0d3c1f3f6c3ec4f66.de-de.xlf"
Add-Link $ws3 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cb5143149f701958c0c10bafc27c1c4e8285a71e/e2e/e74c1a4c-b419-43a7-8bd7-7f98d9726133.md" "e74c1a4c-b419-43a7-8bd7-7f98d9726133.md"
Add-Link $ws3 "C5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97007c52653b3a1097edffc98ab2dcd10ec418bb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/e74c1a4c-b419-43a7-8bd7-7f98d9726133.c2bc1a05f8e55b2535db8e200e39f576949ee6dd.de-de.xlf" "e74c1a4c-b419-43a7-8bd7-7f98d9726133.c2bc1a05f8e55b2535db8e200e39f576949ee6dd.de-de.xlf"
Add-Link $ws3 "A6" "https://github.com/OpenLocalizationTest/oltest/blob/0a7d79782dc0e4eaac7da25566514f6648f5c4ee/e2e/0034157f-dd4d-417a-9df6-110b159bec4f.md" "0034157f-dd4d-417a-9df6-110b159bec4f.md"
Add-Link $ws3 "C6" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/143bf1c2d933fc12d779a0316db24470734d4fae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/0034157f-dd4d-417a-9df6-110b159bec4f.e05f1b9a6487db1a122b3f7a99a8da13ba3f092c.de-de.xlf" "0034157f-dd4d-417a-9df6-110b159bec4f.e05f1b9a6487db1a122b3f7a99a8da13ba3f092c.de-de.xlf"
Add-Link $ws3 "A7" "https://github.com/OpenLocalizationTest/oltest/blob/802a9e18906281c65bb2f5841b7dceee7a668b27/e2e/2291a6aa-fa1a-4622-b6e6-71faa470e61f.md" $newFileMd
Add-Link $ws3 "C7" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19063947110062f2263e7832601bf1ab24f4d384/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/2291a6aa-fa1a-4622-b6e6-71faa470e61f.ee29a6e0dc16ac4a8338fd020a8937d60ebc5a81.de-de.xlf" $newFileDe
Add-Link $ws3 "A8" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/e2e/27a3f539-08b5-4bbe-9811-281f0ccb21ed.md" "27a3f539-08b5-4bbe-9811-281f0ccb21ed.md"
Add-Link $ws3 "C8" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e746485c9f566737a2ce5601f3007deb8795d5c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/ht/27a3f539-08b5-4bbe-9811-281f0ccb21ed.87a85293aeafa20d19b587592fc53056c03f5bcf.de-de.xlf" "27a3f539-08b5-4bbe-9811-281f0ccb21ed.87a85293aeafa20d19b587592fc53056c03f5bcf.de-de.xlf"
Add-Link $ws3 "A9" "https://github.com/OpenLocalizationTest/oltest/blob/e074782d51fb0daf54aa3a76ef3dd7a65de6b9c2/.localization-config" ".localization-config"

Write-Output "done"
